$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("P6").Value = 3.98

# Row 8
$ws.Range("G8").Value = 1.91
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 13
$ws.Range("O8").Value = 1.22
$ws.Range("P8").Value = 4
$ws.Range("Q8").Value = 1.75
$ws.Range("R8").Value = 2.05
$ws.Range("Y8").Value = 8.5
$ws.Range("AC8").Value = 12
$ws.Range("AO8").Value = 10
$ws.Range("AP8").Value = 19

# Row 9
$ws.Range("H9").Value = 3.8
$ws.Range("K9").Value = 2.3
$ws.Range("L9").Value = 5.5
$ws.Range("M9").Value = 1.04
$ws.Range("N9").Value = 12
$ws.Range("O9").Value = 1.22
$ws.Range("P9").Value = 4
$ws.Range("Q9").Value = 1.75
$ws.Range("R9").Value = 2.05
$ws.Range("S9").Value = 1.33
$ws.Range("T9").Value = 3.25
$ws.Range("U9").Value = 1.8
$ws.Range("V9").Value = 1.91
$ws.Range("W9").Value = 7.5
$ws.Range("X9").Value = 8
$ws.Range("AA9").Value = 13
$ws.Range("AB9").Value = 23
$ws.Range("AC9").Value = 12
$ws.Range("AF9").Value = 51
$ws.Range("AG9").Value = 201
$ws.Range("AH9").Value = 15
$ws.Range("AN9").Value = 3.6
$ws.Range("AP9").Value = 19
$ws.Range("AS9").Value = 126
$ws.Range("AT9").Value = 3.25
$ws.Range("AX9").Value = 29
$ws.Range("AY9").Value = 34
$ws.Range("AZ9").Value = 101
$ws.Range("BB9").Value = 201

# Row 10
$ws.Range("G10").Value = 1.8
$ws.Range("I10").Value = 4.5
$ws.Range("J10").Value = 2.5
$ws.Range("L10").Value = 5
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 8
$ws.Range("W10").Value = 6
$ws.Range("X10").Value = 7.5
$ws.Range("AB10").Value = 34
$ws.Range("AC10").Value = 8
$ws.Range("AF10").Value = 67
$ws.Range("AH10").Value = 11
$ws.Range("AI10").Value = 23
$ws.Range("AK10").Value = 51
$ws.Range("AL10").Value = 41
$ws.Range("AO10").Value = 9.5
$ws.Range("AU10").Value = 9
$ws.Range("AX10").Value = 26
$ws.Range("AZ10").Value = 101
$ws.Range("BA10").Value = 126
$ws.Range("BB10").Value = 301

# Row 11
$ws.Range("G11").Value = 2.05
$ws.Range("I11").Value = 3.7
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 9
$ws.Range("O11").Value = 1.4
$ws.Range("P11").Value = 2.75
$ws.Range("Q11").Value = 2.25
$ws.Range("R11").Value = 1.62
$ws.Range("AH11").Value = 8.5
$ws.Range("AN11").Value = 4
$ws.Range("AO11").Value = 12

# Row 12
$ws.Range("G12").Value = 2.35
$ws.Range("I12").Value = 3.1
$ws.Range("J12").Value = 3.2
$ws.Range("N12").Value = 7.5
$ws.Range("W12").Value = 6.5
$ws.Range("X12").Value = 10
$ws.Range("AC12").Value = 7.5
$ws.Range("AG12").Value = 1250
$ws.Range("AI12").Value = 13
$ws.Range("AN12").Value = 4.33
$ws.Range("AO12").Value = 15
$ws.Range("AX12").Value = 19

# Row 13
$ws.Range("G13").Value = 1.4
$ws.Range("H13").Value = 4.2
$ws.Range("I13").Value = 8.5
$ws.Range("L13").Value = 7.5
$ws.Range("O13").Value = 1.29
$ws.Range("P13").Value = 3.5
$ws.Range("Q13").Value = 1.93
$ws.Range("R13").Value = 1.93
$ws.Range("U13").Value = 2.1
$ws.Range("V13").Value = 1.67
$ws.Range("W13").Value = 6.5
$ws.Range("X13").Value = 6.5
$ws.Range("Z13").Value = 9
$ws.Range("AB13").Value = 29
$ws.Range("AD13").Value = 8
$ws.Range("AM13").Value = 51
$ws.Range("AN13").Value = 3.25
$ws.Range("AO13").Value = 7
$ws.Range("AQ13").Value = 21
$ws.Range("AU13").Value = 9.5
$ws.Range("AW13").Value = 8.5
$ws.Range("AZ13").Value = 151
